$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '66.994.48'
Set-TextValue $ws.Range('E2') '  +2.80%  '
Set-TextValue $ws.Range('D3') '3.847.59'
Set-TextValue $ws.Range('E3') '  +4.93%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  -0.16%  '
Set-TextValue $ws.Range('D5') '423.24'
Set-TextValue $ws.Range('E5') '  +3.74%  '
Set-TextValue $ws.Range('D6') '129.34'
Set-TextValue $ws.Range('E6') '  -3.67%  '
Set-TextValue $ws.Range('D7') '3.836.71'
Set-TextValue $ws.Range('E7') '  +4.78%  '
Set-TextValue $ws.Range('D8') '0.608'
Set-TextValue $ws.Range('E8') '  -2.31%  '
Set-TextValue $ws.Range('D9') '0.998'
Set-TextValue $ws.Range('E9') '  -0.23%  '
Set-TextValue $ws.Range('D10') '0.722'
Set-TextValue $ws.Range('E10') '  -0.97%  '
Set-TextValue $ws.Range('D11') '0.159'
Set-TextValue $ws.Range('E11') '  -2.65%  '
Set-TextValue $ws.Range('D12') '0.0000345'
Set-TextValue $ws.Range('E12') '  +6.14%  '
Set-TextValue $ws.Range('D13') '40.81'
Set-TextValue $ws.Range('E13') '  -3.69%  '
Set-TextValue $ws.Range('B14') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D14') '4.454.81'
Set-TextValue $ws.Range('E14') '  +4.31%  '
Set-TextValue $ws.Range('B15') 'Polkadot'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D15') '10.22'
Set-TextValue $ws.Range('E15') '  +2.33%  '
Set-TextValue $ws.Range('D16') '15.77'
Set-TextValue $ws.Range('E16') '  +14.45%  '
Set-TextValue $ws.Range('B17') 'WrappedEther'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D17') '3.876.36'
Set-TextValue $ws.Range('E17') '  +5.21%  '
Set-TextValue $ws.Range('B18') 'TRON'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D18') '0.138'
Set-TextValue $ws.Range('E18') '  -0.64%  '
Set-TextValue $ws.Range('D19') '19.74'
Set-TextValue $ws.Range('E19') '  -1.68%  '
Set-TextValue $ws.Range('D20') '67.246.71'
Set-TextValue $ws.Range('E20') '  +2.85%  '
Set-TextValue $ws.Range('E21') '  -0.84%  '
Set-TextValue $ws.Range('D22') '409.17'
Set-TextValue $ws.Range('E22') '  -3.18%  '
Set-TextValue $ws.Range('D23') '15.02'
Set-TextValue $ws.Range('E23') '  -2.21%  '
Set-TextValue $ws.Range('D24') '84.14'
Set-TextValue $ws.Range('E24') '  -2.40%  '
Set-TextValue $ws.Range('E25') '  +1.42%  '
Set-TextValue $ws.Range('D26') '37.31'
Set-TextValue $ws.Range('E26') '  +3.91%  '
Set-TextValue $ws.Range('E27') '  +6.80%  '
Set-TextValue $ws.Range('D28') '3.23'
Set-TextValue $ws.Range('E28') '  +0.64%  '
Set-TextValue $ws.Range('D29') '5.43'
Set-TextValue $ws.Range('E29') '  +5.58%  '
Set-TextValue $ws.Range('D30') '9.32'
Set-TextValue $ws.Range('E30') '  +33.21%  '
Set-TextValue $ws.Range('D31') '742.46'
Set-TextValue $ws.Range('E31') '  +8.83%  '
Set-TextValue $ws.Range('D32') '13.03'
Set-TextValue $ws.Range('E32') '  +2.22%  '
Set-TextValue $ws.Range('E33') '  +2.25%  '
Set-TextValue $ws.Range('E34') '  +2.39%  '
Set-TextValue $ws.Range('D35') '0.999'
Set-TextValue $ws.Range('E35') '  -0.13%  '
Set-TextValue $ws.Range('E36') '  -6.54%  '
Set-TextValue $ws.Range('D37') '38.56'
Set-TextValue $ws.Range('E37') '  -7.40%  '
Set-TextValue $ws.Range('D38') '55.85'
Set-TextValue $ws.Range('E38') '  -0.23%  '
Set-TextValue $ws.Range('D39') '5.41'
Set-TextValue $ws.Range('E39') '  +22.78%  '
Set-TextValue $ws.Range('D40') '0.0₃0751'
Set-TextValue $ws.Range('E40') '  +15.32%  '
Set-TextValue $ws.Range('E41') '  -2.47%  '
Set-TextValue $ws.Range('D42') '2.89'
Set-TextValue $ws.Range('E42') '  -1.87%  '
Set-TextValue $ws.Range('E43') '  +0.58%  '
Set-TextValue $ws.Range('D44') '3.36'
Set-TextValue $ws.Range('E44') '  +0.37%  '
Set-TextValue $ws.Range('E45') '  -4.64%  '
Set-TextValue $ws.Range('D46') '3.12'
Set-TextValue $ws.Range('E46') '  +0.04%  '
Set-TextValue $ws.Range('B47') 'TheGraph'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D47') '0.313'
Set-TextValue $ws.Range('E47') '  +7.24%  '
Set-TextValue $ws.Range('B48') 'Monero'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D48') '142.80'
Set-TextValue $ws.Range('E48') '  -0.56%  '
Set-TextValue $ws.Range('D49') '2.05'
Set-TextValue $ws.Range('E49') '  -2.04%  '
Set-TextValue $ws.Range('D50') '2.82'
Set-TextValue $ws.Range('E50') '  +0.39%  '
Set-TextValue $ws.Range('B51') 'WEMIXToken'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D51') '2.55'
Set-TextValue $ws.Range('E51') '  +0.73%  '
